$d = $word.ActiveDocument

function Find-Replace($findText, $replaceText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        Write-Output "NOT FOUND (replace): $findText"
    }
    return $found
}

function Delete-ParaText($findText) {
    # Finds findText (which should span an entire paragraph's visible text),
    # and deletes it together with its trailing paragraph mark.
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND (delete): $findText"
        return
    }
    $delRng = $d.Range($rng.Start, $rng.End + 1)
    $delRng.Delete()
}

function Insert-ParaBefore($findText, $newParaText) {
    # Finds findText and inserts a new paragraph containing newParaText
    # immediately before the paragraph that findText starts.
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND (insert-before): $findText"
        return
    }
    $ins = $d.Range($rng.Start, $rng.Start)
    $ins.InsertBefore($newParaText + "`r")
}

# 1) "Magic Elements are Fire-Ice-Lightning-Null-Light." -> "...Divine-Fire-Ice-Lightning."
Find-Replace "Fire-Ice-Lightning-Null-Light" "Divine-Fire-Ice-Lightning"

# 2) Basic/Advanced/Expert spell descriptions: Enemy/enemies -> target/targets, hits -> hit
Find-Replace "Basic Spells hit only one Enemy." "Basic Spells only hit one target."
Find-Replace "Advanced Spells hits multiple enemies, via mini-radius" "Advanced Spells hit multiple targets, via mini-radius"
Find-Replace "Expert Spells hits multiple enemies with a much greater radius" "Expert Spells hit multiple targets with a much greater radius"

# 3) Insert new "Divine" spell-tree paragraph before the "Fire" paragraph.
Insert-ParaBefore "Fire: Fireball -> Fire Blast -> Rain of Fire" "Divine: Divine Heal -> Divine Barrier -> Sanctuary"

# 4) Fire paragraph: "Rain of Fire" -> "Armageddon"
Find-Replace "Fire: Fireball -> Fire Blast -> Rain of Fire" "Fire: Fireball -> Fire Blast -> Armageddon"

# 5) Remove the old "Thunder: ..." paragraph and replace it with a new
#    "Lightning: ..." paragraph (placed, for now, right where Thunder was;
#    it will end up after Ice once Ice is relocated below).
Delete-ParaText "Thunder: Lightning Bolt -> Chain Lightning -> ThunderStorm"
Insert-ParaBefore "Ice: Ice Shard -> Mist -> Blizzard" "Lightning: Lightning Bolt -> Chain Lightning -> Thunderstorm"

# 6) Ice paragraph: "Mist" -> "Ice Frost", then move it above the new Lightning paragraph.
Find-Replace "Ice: Ice Shard -> Mist -> Blizzard" "Ice: Ice Shard -> Ice Frost -> Blizzard"
Delete-ParaText "Ice: Ice Shard -> Ice Frost -> Blizzard"
Insert-ParaBefore "Lightning: Lightning Bolt -> Chain Lightning -> Thunderstorm" "Ice: Ice Shard -> Ice Frost -> Blizzard"

# 7) Remove the "Light: ..." and "Null: ..." paragraphs entirely.
Delete-ParaText "Light: Heal -> Heal Barrier -> Heal All"
Delete-ParaText "Null: Magic Missile-> Missile Barrage -> Magic Detonation"

Write-Output "done"
